# Fluxo de cadastro de gastos com parcelamento real com uma linha por parcela
#
# 1. GASTOS header row: rename the "parcelas" columns so each installment
#    gets its own row (Parcelado / Numero  Parcelas / Valor Parcela).
# 2. ENTRADAS: a new (empty but underline-formatted) row is appended below
#    the existing blank row, extending the used range to G6, and the
#    selection moves to H1.
# 3. GASTOS becomes the active sheet (selection A2), ENTRADAS loses the
#    active/tabSelected flag.

$wb = $excel.ActiveWorkbook

$entradas = $wb.Worksheets.Item("ENTRADAS")
$gastos   = $wb.Worksheets.Item("GASTOS")

# --- GASTOS: rename the installment-related headers -----------------------
$gastos.Range("H1").Value = "Parcelado"
$gastos.Range("I1").Value = "Numero  Parcelas"
$gastos.Range("J1").Value = "Valor Parcela"

# --- ENTRADAS: add a new formatted (underlined) blank row below row 5 -----
$entradas.Rows.Item(6).Font.Underline = $true

# --- Selections -------------------------------------------------------------
$entradas.Range("H1").Select()

# --- Make GASTOS the active / selected sheet, matching the saved view -----
$gastos.Activate()
$gastos.Range("A2").Select()

Write-Host "done"
